# Add a new "2022-Q4" quarterly sheet, positioned right after "总计" and
# before the existing "2022-Q2" sheet, and refresh the "总计" (totals)
# summary sheet with the new quarter's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by copying the existing "2022-Q2"
#    sheet (same column layout/styling), placed immediately before it.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($templateSheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q4"

# The template only has 2 data rows (rows 2-3); we need 5, so clone the
# formatting of row 2 down into rows 4-6 before writing values.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A4:H4").PasteSpecial(-4122)
$newSheet.Range("A5:H5").PasteSpecial(-4122)
$newSheet.Range("A6:H6").PasteSpecial(-4122)

# Make sure the text-like columns (fund code, name, size, position,
# ratio, market value) stay as text instead of being auto-coerced into
# numbers/dates by the smart-entry logic.
$newSheet.Range("B2:G6").NumberFormat = "@"

$data = @(
    @(0, "007832", "博道伍佰智航股票C",       "6.00", "93.04", "0.90", "0.0540", 2),
    @(1, "013466", "博时智选量化多因子股票C", "2.26", "93.55", "1.54", "0.0348", 2),
    @(2, "007831", "博道伍佰智航股票A",       "2.75", "93.04", "0.90", "0.0248", 7),
    @(3, "003238", "新华外延增长主题灵活配置混合", "0.51", "85.17", "3.23", "0.0165", 2),
    @(4, "013465", "博时智选量化多因子股票A", "0.67", "93.55", "1.54", "0.0103", 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $vals[0]
    $newSheet.Cells.Item($row, 2).Value = $vals[1]
    $newSheet.Cells.Item($row, 3).Value = $vals[2]
    $newSheet.Cells.Item($row, 4).Value = $vals[3]
    $newSheet.Cells.Item($row, 5).Value = $vals[4]
    $newSheet.Cells.Item($row, 6).Value = $vals[5]
    $newSheet.Cells.Item($row, 7).Value = $vals[6]
    $newSheet.Cells.Item($row, 8).Value = $vals[7]
}

# ---------------------------------------------------------------------
# 2) Refresh the "总计" (totals) sheet: insert the new 2022-Q4 row at
#    the top of the data and push the existing quarters down, adding
#    the previously-missing 2021-Q4 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Give the new row 5 the same index-column styling as the rows above it.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalData = @(
    @(0, "2022-Q4", 5, 0.14),
    @(1, "2022-Q2", 2, 0),
    @(2, "2021-Q4", 2, 0.04),
    @(3, "2021-Q1", 1, 0.04)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $row = $i + 2
    $vals = $totalData[$i]
    $totalSheet.Cells.Item($row, 1).Value = $vals[0]
    $totalSheet.Cells.Item($row, 2).Value = $vals[1]
    $totalSheet.Cells.Item($row, 3).Value = $vals[2]
    $totalSheet.Cells.Item($row, 4).Value = $vals[3]
}

# Keep "总计" as the active/selected tab (unchanged from the original
# workbook) rather than leaving the newly-added sheet activated.
$totalSheet.Activate()
